$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new defined name referencing the new uncertainty cell
$wb.Names.Add("UNC_Batt_Cap", "=Sheet1!`$B`$9")

# Update the Battery_Capacity formula to include the new uncertainty factor
$ws.Range("B5").Formula = "=Pack_Count*5*11.1*3600*UNC_Batt_Cap"

# New row 9: label, value, unit
$ws.Range("A9").Value = "UNC_Batt_Cap"
$ws.Range("B9").Value = 0.76
$ws.Range("C9").Value = "% Battery Carged"

# Restore the active selection to B6 as in the target workbook
$selResult = $ws.Range("B6").Select()

$saveResult = $wb.Save()

Write-Host "done"
